# Updated symbol list on Tue Dec 20 21:51:49 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Price-only updates (column D) ---
Set-TextValue "D2"  "250.66"
Set-TextValue "D3"  "22.88"
Set-TextValue "D4"  "5.426"
Set-TextValue "D7"  "6.372"
Set-TextValue "D8"  "0.8137"
Set-TextValue "D9"  "0.9231"

# --- Rows 10-18: coin list shifted up by one position, new coin appended at the end ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1442"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07475"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03128"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03059"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09353"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.723"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001594"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04776"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005791"
$ws.Range("E18").Value = "17OneONE"

# --- More price-only updates (column D) ---
Set-TextValue "D19" "0.006380"
Set-TextValue "D20" "0.005040"
Set-TextValue "D24" "2.178"
Set-TextValue "D28" "0.0003000"
Set-TextValue "D40" "0.04030"

# --- Row 41: price + label change ---
Set-TextValue "D41" "0.002940"
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"

# --- More price-only updates (column D) ---
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.002710"
Set-TextValue "D44" "0.008018"
Set-TextValue "D45" "0.00005803"

# --- Label-only updates (column E) ---
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
